$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.929.97'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.811.92'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.42'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.58%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.01%  '
$ws.Range('E7').ClearFormats()
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3893'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.89%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09467'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +19.55%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.100'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.07%  '
$ws.Range('E10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.43'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.23%  '
$ws.Range('E11').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.415'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.06%  '
$ws.Range('E12').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.52'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.813.52'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.30%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.284'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('E16').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.69%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.51'
$ws.Range('D18').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06594'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('E20').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.16'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.71%  '
$ws.Range('E21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.953'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('E22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.985.83'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.74%  '
$ws.Range('E23').ClearFormats()
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.20'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('E24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.248'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('E25').ClearFormats()
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Monero'
$ws.Range('B26').ClearFormats()
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C26').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.78'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('E26').ClearFormats()
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('B27').ClearFormats()
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C27').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.75'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.46%  '
$ws.Range('E27').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.020.72'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.01%  '
$ws.Range('E28').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.401'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.73%  '
$ws.Range('E29').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.69'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.59%  '
$ws.Range('E30').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.81%  '
$ws.Range('E31').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.056'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('E32').ClearFormats()
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.628'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.05%  '
$ws.Range('E33').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.627'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.30%  '
$ws.Range('E34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.06825'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -6.55%  '
$ws.Range('E35').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.012'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.87%  '
$ws.Range('E36').ClearFormats()
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02320'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E37').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2151'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.43%  '
$ws.Range('E38').ClearFormats()
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -7.07%  '
$ws.Range('E39').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.951'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.31%  '
$ws.Range('E40').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6253'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.58%  '
$ws.Range('E41').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9998'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('E42').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.146'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.22%  '
$ws.Range('E43').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.08'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('E44').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5918'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.89%  '
$ws.Range('E45').ClearFormats()
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.293'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -5.68%  '
$ws.Range('E46').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.692'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.41%  '
$ws.Range('E47').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.23'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.77%  '
$ws.Range('E48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.961'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.23%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.180'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.32%  '
$ws.Range('E50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06776'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.11%  '
$ws.Range('E51').ClearFormats()
